$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.144.56"
$ws.Range("E2").Value = "  +1.19%  "

$ws.Range("D3").Value = "1.912.92"
$ws.Range("E3").Value = "  +1.77%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'327.92"
$ws.Range("E5").Value = "  +1.08%  "

$ws.Range("D6").Value = "'1.006"
$ws.Range("E6").Value = "  +0.01%  "

$ws.Range("D7").Value = "'0.4660"
$ws.Range("E7").Value = "  -0.37%  "

$ws.Range("D8").Value = "'0.3956"
$ws.Range("E8").Value = "  +0.50%  "

$ws.Range("D9").Value = "'47.18"
$ws.Range("E9").Value = "  +1.22%  "

$ws.Range("D10").Value = "'0.08016"
$ws.Range("E10").Value = "  +1.17%  "

$ws.Range("D11").Value = "'1.012"
$ws.Range("E11").Value = "  +3.08%  "

$ws.Range("D12").Value = "'22.33"
$ws.Range("E12").Value = "  -0.22%  "

$ws.Range("D13").Value = "1.913.97"
$ws.Range("E13").Value = "  -0.66%  "

$ws.Range("D14").Value = "'7.153"
$ws.Range("E14").Value = "  +1.92%  "

$ws.Range("D15").Value = "'5.805"
$ws.Range("E15").Value = "  +1.09%  "

$ws.Range("D16").Value = "'0.06962"
$ws.Range("E16").Value = "  -0.28%  "

$ws.Range("D17").Value = "'89.31"
$ws.Range("E17").Value = "  +0.65%  "

$ws.Range("E18").Value = "  +0.05%  "

$ws.Range("D19").Value = "'0.00001020"
$ws.Range("E19").Value = "  +0.99%  "

$ws.Range("D20").Value = "'17.32"
$ws.Range("E20").Value = "  +2.09%  "

$ws.Range("E21").Value = "  -0.08%  "

$ws.Range("D22").Value = "29.153.68"
$ws.Range("E22").Value = "  +1.14%  "

$ws.Range("D23").Value = "'5.415"
$ws.Range("E23").Value = "  +1.22%  "

$ws.Range("D24").Value = "'11.16"
$ws.Range("E24").Value = "  +0.55%  "

$ws.Range("D25").Value = "2.154.88"
$ws.Range("E25").Value = "  +1.73%  "

$ws.Range("E26").Value = "  -3.44%  "

$ws.Range("D27").Value = "'156.00"
$ws.Range("E27").Value = "  +1.64%  "

$ws.Range("D28").Value = "'19.64"
$ws.Range("E28").Value = "  +1.19%  "

$ws.Range("D29").Value = "'5.942"
$ws.Range("E29").Value = "  +2.79%  "

$ws.Range("D30").Value = "'2.032"
$ws.Range("E30").Value = "  +1.54%  "

$ws.Range("D31").Value = "'120.83"
$ws.Range("E31").Value = "  +0.78%  "

$ws.Range("D32").Value = "'0.09402"
$ws.Range("E32").Value = "  +0.01%  "

$ws.Range("D33").Value = "'0.9429"
$ws.Range("E33").Value = "  +0.63%  "

$ws.Range("D34").Value = "'5.376"
$ws.Range("E34").Value = "  +1.21%  "

$ws.Range("D35").Value = "'1.348"
$ws.Range("E35").Value = "  -0.72%  "

$ws.Range("D36").Value = "'3.283"
$ws.Range("E36").Value = "  -1.96%  "

$ws.Range("D37").Value = "'0.05886"
$ws.Range("E37").Value = "  -0.43%  "

$ws.Range("E38").Value = "  +3.37%  "

$ws.Range("E39").Value = "  +0.75%  "

$ws.Range("D40").Value = "'0.02105"
$ws.Range("E40").Value = "  -1.16%  "

$ws.Range("D41").Value = "'0.5848"
$ws.Range("E41").Value = "  +2.06%  "

$ws.Range("D42").Value = "'1.006"
$ws.Range("E42").Value = "  +0.10%  "

$ws.Range("D43").Value = "'0.1822"
$ws.Range("E43").Value = "  +1.35%  "

$ws.Range("D44").Value = "'10.06"
$ws.Range("E44").Value = "  +0.55%  "

$ws.Range("D45").Value = "'2.300"
$ws.Range("E45").Value = "  +10.56%  "

$ws.Range("D46").Value = "'12.04"
$ws.Range("E46").Value = "  +1.59%  "

$ws.Range("D47").Value = "'0.5478"
$ws.Range("E47").Value = "  +2.24%  "

$ws.Range("D48").Value = "'0.07207"
$ws.Range("E48").Value = "  -1.43%  "

$ws.Range("D49").Value = "'1.915"
$ws.Range("E49").Value = "  +3.65%  "

$ws.Range("D50").Value = "'1.130"
$ws.Range("E50").Value = "  -4.31%  "

$ws.Range("E51").Value = "  -0.81%  "
